$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '22.548.46'
$ws.Range("E2").Value = '  +0.29%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.577.01'
$ws.Range("E3").Value = '  +0.33%  '
$ws.Range("E4").Value = '  +0.03%  '
$ws.Range("E5").Value = '  +0.01%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '289.30'
$ws.Range("E6").Value = '  -0.48%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.3733'
$ws.Range("E7").Value = '  +0.66%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '48.36'
$ws.Range("E8").Value = '  -3.36%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.3356'
$ws.Range("E9").Value = '  -0.60%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '1.139'
$ws.Range("E10").Value = '  -1.06%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07517'
$ws.Range("E11").Value = '  -0.33%  '
$ws.Range("E12").Value = '  +0.02%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '21.05'
$ws.Range("E13").Value = '  -0.66%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '5.998'
$ws.Range("E14").Value = '  -0.40%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '6.954'
$ws.Range("E15").Value = '  -0.26%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '1.583.20'
$ws.Range("E16").Value = '  +0.73%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.00001123'
$ws.Range("E17").Value = '  +0.12%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '88.75'
$ws.Range("E18").Value = '  -1.91%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.06776'
$ws.Range("E19").Value = '  +0.11%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '6.413'
$ws.Range("E20").Value = '  +0.95%  '
$ws.Range("E21").Value = '  -0.01%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '16.56'
$ws.Range("E22").Value = '  +0.79%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '12.15'
$ws.Range("E23").Value = '  -0.74%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '22.547.43'
$ws.Range("E24").Value = '  +0.33%  '
$ws.Range("E25").Value = '  +1.05%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '2.591'
$ws.Range("E26").Value = '  -1.63%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '153.15'
$ws.Range("E27").Value = '  +2.62%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '19.78'
$ws.Range("E28").Value = '  -1.38%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '5.020'
$ws.Range("E29").Value = '  -1.30%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '124.44'
$ws.Range("E30").Value = '  -0.60%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.757.81'
$ws.Range("E31").Value = '  +0.47%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '1.054'
$ws.Range("E32").Value = '  -1.55%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '6.201'
$ws.Range("E33").Value = '  -0.08%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '2.017'
$ws.Range("E34").Value = '  +0.17%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '9.765'
$ws.Range("E35").Value = '  -0.39%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.08338'
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.02466'
$ws.Range("E37").Value = '  -0.65%  '
$ws.Range("E38").Value = '  -0.48%  '
$ws.Range("B39").Value = 'Hedera'
$ws.Range("C39").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.06414'
$ws.Range("E39").Value = '  -2.02%  '
$ws.Range("B40").Value = 'InternetComputer(DFINITY)'
$ws.Range("C40").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '5.422'
$ws.Range("E40").Value = '  -0.28%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '1.300'
$ws.Range("B42").Value = 'Aptos'
$ws.Range("C42").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '11.37'
$ws.Range("E42").Value = '  +0.59%  '
$ws.Range("B43").Value = 'TheSandbox'
$ws.Range("C43").Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.6337'
$ws.Range("E43").Value = '  +1.72%  '
$ws.Range("E44").Value = '  -0.02%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '13.94'
$ws.Range("E45").Value = '  -1.35%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.6158'
$ws.Range("E46").Value = '  +5.03%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '3.791'
$ws.Range("E47").Value = '  -0.37%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '2.069'
$ws.Range("E48").Value = '  -0.33%  '
$ws.Range("E49").Value = '  -2.71%  '
$ws.Range("E50").Value = '  -1.44%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.07286'
$ws.Range("E51").Value = '  -0.40%  '
